# Add movie scraping functionality and update report with new data
#
# Refreshes the movie report table (A:D = Title, Genre, Rating, Year) on the
# active sheet with freshly scraped values.
#
# Rating/Year are stored as TEXT in this sheet (same as the rest of the
# table), so numeric-looking values are written with a leading apostrophe
# to force Excel to keep them as text instead of silently re-typing the
# cell as a number. Title/Genre are plain (non-numeric) text and are
# written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cell, $value) {
    $ws.Range($cell).Value = $value
}

function Set-NumericText($cell, $value) {
    # Leading apostrophe = Excel "text" quote-prefix; keeps e.g. "5.9"/"2024"
    # stored as text rather than being coerced into a number.
    $ws.Range($cell).Value = "'" + $value
}

# Row 2
Set-Text "A2" "Don't Move"
Set-Text "B2" "Horror"
Set-NumericText "C2" "5.9"

# Row 3
Set-Text "A3" "Joker: Folie à Deux"
Set-NumericText "C3" "5.3"

# Row 4
Set-Text "A4" "Terrifier 3"
Set-NumericText "C4" "6.7"

# Row 5
Set-Text "A5" "Am I Racist?"
Set-Text "B5" "Comedy"
Set-NumericText "C5" "7"

# Row 6
Set-Text "A6" "[FR] Family Pack"
Set-Text "B6" "Adventure"
Set-NumericText "C6" "5.4"

# Row 7
Set-Text "A7" "Transformers One"
Set-Text "B7" "Action"
Set-NumericText "C7" "7.7"

# Row 8
Set-Text "A8" "Carved"
Set-Text "B8" "Comedy"
Set-NumericText "C8" "4.4"

# Row 9
Set-Text "A9" "Bagman"
Set-Text "B9" "Horror"
Set-NumericText "C9" "4.6"

# Row 10
Set-Text "A10" "Die Alone"
Set-Text "B10" "Horror"
Set-NumericText "C10" "5.7"

# Row 11
Set-Text "A11" "Sweet Bobby: My Catfish Nightmare"
Set-Text "B11" "Documentary"
Set-NumericText "C11" "6.2"
Set-NumericText "D11" "2024"

# Row 12
Set-Text "A12" "Death Becomes Her"
Set-Text "B12" "Action"
Set-NumericText "C12" "6.7"
Set-NumericText "D12" "1992"

# Row 13
Set-Text "A13" "Armageddon Time"
Set-Text "B13" "Action"
Set-NumericText "C13" "6.5"
Set-NumericText "D13" "2022"

# Row 14
Set-Text "A14" "Cowboys & Aliens"
Set-NumericText "C14" "6"
Set-NumericText "D14" "2011"

# Row 15
Set-Text "A15" "Oddity"
Set-Text "B15" "Horror"
Set-NumericText "C15" "6.7"
Set-NumericText "D15" "2024"

# Row 16
Set-Text "A16" "Cuckoo"
Set-Text "B16" "Horror"
Set-NumericText "D16" "2024"

# Row 17
Set-Text "A17" "Friday the 13th"
Set-NumericText "C17" "5.5"
Set-NumericText "D17" "2009"

# Row 18
Set-Text "A18" "Slingshot"
Set-Text "B18" "Sci-Fi"
Set-NumericText "C18" "5.7"
Set-NumericText "D18" "2024"

# Row 19
Set-Text "A19" "A Nightmare on Elm Street"
Set-NumericText "C19" "7.4"
Set-NumericText "D19" "1984"

# Row 20
Set-Text "A20" "You're Next"
Set-NumericText "C20" "6.6"
Set-NumericText "D20" "2011"

# Row 21
Set-Text "A21" "Night Train"
Set-NumericText "C21" "3.7"
Set-NumericText "D21" "2023"
